$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 51-52 (existing row 51 shifts down to row 53)
$ws.Range("A51:A52").EntireRow.Insert()

# New row 51
$ws.Range("A51").Value = 7
$ws.Range("B51").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C51").Value = "Ñuble"
$ws.Range("D51").Value = 45265
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100101
$ws.Range("H51").Value = "Berries"
$ws.Range("I51").Value = 100101001
$ws.Range("J51").Value = "Arándano (blue)"
$ws.Range("K51").Value = "Sin especificar"
$ws.Range("L51").Value = "Primera"
$ws.Range("M51").Value = 100
$ws.Range("N51").Value = 6000
$ws.Range("O51").Value = 6000
$ws.Range("P51").Value = 6000
$ws.Range("Q51").Value = "$/bandeja 2 kilos"
$ws.Range("R51").Value = "Provincia de Diguillín"
$ws.Range("S51").Value = 3000
$ws.Range("T51").Value = 2

# New row 52
$ws.Range("A52").Value = 7
$ws.Range("B52").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C52").Value = "Ñuble"
$ws.Range("D52").Value = 45265
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100101
$ws.Range("H52").Value = "Berries"
$ws.Range("I52").Value = 100101001
$ws.Range("J52").Value = "Arándano (blue)"
$ws.Range("K52").Value = "Sin especificar"
$ws.Range("L52").Value = "Segunda"
$ws.Range("M52").Value = 100
$ws.Range("N52").Value = 5000
$ws.Range("O52").Value = 5000
$ws.Range("P52").Value = 5000
$ws.Range("Q52").Value = "$/bandeja 2 kilos"
$ws.Range("R52").Value = "Provincia de Diguillín"
$ws.Range("S52").Value = 2500
$ws.Range("T52").Value = 2

# Ensure the date cells use the existing date number format (style index 2 from styles.xml)
$ws.Range("D51").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D52").NumberFormat = "YYYY-MM-DD HH:MM:SS"
